$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# Add header labels for the newly inserted row
$ws.Range("A1").Value = "state"
$ws.Range("B1").Value = "party"

# Update the selection to match the final state
$ws.Range("D6").Select()
